# Auto-generated Excel COM-interop script to update Kujata_Profits leve profit data
# across all affected sheets, matching the scheduled-runner market-price refresh diff.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 499.08163
$ws.Range("J17").Value = 499.08163
$ws.Range("L17").Value = 1497.24489
$ws.Range("N17").Value = -1833.24489
$ws.Range("H33").Value = 727.8889
$ws.Range("I33").Value = 865.1667
$ws.Range("J33").Value = 453.33334
$ws.Range("K33").Value = 865.1667
$ws.Range("L33").Value = 453.33334
$ws.Range("M33").Value = -636.1667
$ws.Range("N33").Value = -911.33334
$ws.Range("H98").Value = 13054.615
$ws.Range("I98").Value = 6337.273
$ws.Range("K98").Value = 6337.273
$ws.Range("M98").Value = -4839.273
$ws.Range("H115").Value = 658.3333
$ws.Range("I115").Value = 658.3333
$ws.Range("K115").Value = 1974.9999
$ws.Range("M115").Value = -407.9999
$ws.Range("H122").Value = 13054.615
$ws.Range("I122").Value = 6337.273
$ws.Range("K122").Value = 19011.819
$ws.Range("M122").Value = -16561.819
$ws.Range("H131").Value = 795.7143
$ws.Range("I131").Value = 792.5
$ws.Range("J131").Value = 800
$ws.Range("K131").Value = 2377.5
$ws.Range("L131").Value = 2400
$ws.Range("M131").Value = 2662.5
$ws.Range("N131").Value = -12480

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17620.834
$ws.Range("I2").Value = 800
$ws.Range("J2").Value = 26031.25
$ws.Range("K2").Value = 800
$ws.Range("L2").Value = 26031.25
$ws.Range("M2").Value = -687
$ws.Range("N2").Value = -26257.25
$ws.Range("H32").Value = 3149.5
$ws.Range("I32").Value = 2915.2388
$ws.Range("J32").Value = 4576.364
$ws.Range("K32").Value = 2915.2388
$ws.Range("L32").Value = 4576.364
$ws.Range("M32").Value = -2628.2388
$ws.Range("N32").Value = -5150.364
$ws.Range("H42").Value = 5915.5
$ws.Range("I42").Value = 1800
$ws.Range("K42").Value = 1800
$ws.Range("M42").Value = -1314
$ws.Range("H61").Value = 1048.4318
$ws.Range("I61").Value = 884.2368
$ws.Range("K61").Value = 884.2368
$ws.Range("M61").Value = -672.2368
$ws.Range("H74").Value = 1393.091
$ws.Range("I74").Value = 861.3333
$ws.Range("J74").Value = 2532.5715
$ws.Range("K74").Value = 861.3333
$ws.Range("L74").Value = 2532.5715
$ws.Range("M74").Value = 12.66669999999999
$ws.Range("N74").Value = -4280.5715
$ws.Range("H77").Value = 1393.091
$ws.Range("I77").Value = 861.3333
$ws.Range("J77").Value = 2532.5715
$ws.Range("K77").Value = 4306.6665
$ws.Range("L77").Value = 12662.8575
$ws.Range("M77").Value = 61.33349999999973
$ws.Range("N77").Value = -21398.8575
$ws.Range("H116").Value = 17620.834
$ws.Range("I116").Value = 800
$ws.Range("J116").Value = 26031.25
$ws.Range("K116").Value = 800
$ws.Range("L116").Value = 26031.25
$ws.Range("M116").Value = 1494
$ws.Range("N116").Value = -30619.25
$ws.Range("H122").Value = 945.85
$ws.Range("I122").Value = 970.41174
$ws.Range("K122").Value = 2911.23522
$ws.Range("M122").Value = -461.23522
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 1638.75
$ws.Range("I132").Value = 1328.625
$ws.Range("K132").Value = 3985.875
$ws.Range("M132").Value = -1455.875
$ws.Range("H136").Value = 1048.4318
$ws.Range("I136").Value = 884.2368
$ws.Range("K136").Value = 2652.7104
$ws.Range("M136").Value = -102.7103999999999
$ws.Range("H137").Value = 35945
$ws.Range("J137").Value = 35945
$ws.Range("L137").Value = 35945
$ws.Range("N137").Value = -46145
$ws.Range("H139").Value = 37010
$ws.Range("J139").Value = 37010
$ws.Range("L139").Value = 37010
$ws.Range("N139").Value = -47290

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17620.834
$ws.Range("I3").Value = 800
$ws.Range("J3").Value = 26031.25
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 26031.25
$ws.Range("M3").Value = -686
$ws.Range("N3").Value = -26259.25
$ws.Range("H99").Value = 33334516
$ws.Range("I99").Value = 45455630
$ws.Range("J99").Value = 1450
$ws.Range("K99").Value = 45455630
$ws.Range("L99").Value = 1450
$ws.Range("M99").Value = -45454132
$ws.Range("N99").Value = -4446
$ws.Range("H134").Value = 4850.3335
$ws.Range("I134").Value = 1224.8096
$ws.Range("J134").Value = 11195
$ws.Range("K134").Value = 3674.4288
$ws.Range("L134").Value = 33585
$ws.Range("M134").Value = -1139.4288
$ws.Range("N134").Value = -38655

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1213.78
$ws.Range("I31").Value = 1242.5333
$ws.Range("J31").Value = 955
$ws.Range("K31").Value = 1242.5333
$ws.Range("L31").Value = 955
$ws.Range("M31").Value = -947.5333000000001
$ws.Range("N31").Value = -1545
$ws.Range("H34").Value = 1213.78
$ws.Range("I34").Value = 1242.5333
$ws.Range("J34").Value = 955
$ws.Range("K34").Value = 1242.5333
$ws.Range("L34").Value = 955
$ws.Range("M34").Value = -1040.5333
$ws.Range("N34").Value = -1359
$ws.Range("H99").Value = 2393867.2
$ws.Range("I99").Value = 2925393.2
$ws.Range("K99").Value = 2925393.2
$ws.Range("M99").Value = -2923895.2
$ws.Range("H126").Value = 2393867.2
$ws.Range("I126").Value = 2925393.2
$ws.Range("K126").Value = 8776179.600000001
$ws.Range("M126").Value = -8773709.600000001
$ws.Range("H132").Value = 3034.4
$ws.Range("I132").Value = 2451.8
$ws.Range("K132").Value = 7355.400000000001
$ws.Range("M132").Value = -4825.400000000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1196
$ws.Range("J9").Value = 1196
$ws.Range("L9").Value = 3588
$ws.Range("N9").Value = -4036
$ws.Range("H38").Value = 33
$ws.Range("I38").Value = 33
$ws.Range("K38").Value = 99
$ws.Range("M38").Value = 248
$ws.Range("H131").Value = 1432.3
$ws.Range("J131").Value = 1508.9354
$ws.Range("L131").Value = 4526.8062
$ws.Range("N131").Value = -14606.8062

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1470.2307
$ws.Range("I113").Value = 1442.8572
$ws.Range("J113").Value = 1502.1666
$ws.Range("K113").Value = 1442.8572
$ws.Range("L113").Value = 1502.1666
$ws.Range("M113").Value = 727.1428000000001
$ws.Range("N113").Value = -5842.1666
$ws.Range("H122").Value = 1546.5
$ws.Range("I122").Value = 1552.2941
$ws.Range("J122").Value = 1526.8
$ws.Range("K122").Value = 4656.8823
$ws.Range("L122").Value = 4580.4
$ws.Range("M122").Value = -2206.8823
$ws.Range("N122").Value = -9480.4
$ws.Range("H132").Value = 2218.2354
$ws.Range("I132").Value = 1630.0588
$ws.Range("J132").Value = 2806.4119
$ws.Range("K132").Value = 4890.1764
$ws.Range("L132").Value = 8419.235700000001
$ws.Range("M132").Value = -2360.1764
$ws.Range("N132").Value = -13479.2357
$ws.Range("H140").Value = 29488
$ws.Range("J140").Value = 29488
$ws.Range("L140").Value = 29488
$ws.Range("N140").Value = -39848

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1510.0834
$ws.Range("I22").Value = 745
$ws.Range("J22").Value = 1663.1
$ws.Range("K22").Value = 745
$ws.Range("L22").Value = 1663.1
$ws.Range("M22").Value = -450
$ws.Range("N22").Value = -2253.1
$ws.Range("H27").Value = 1510.0834
$ws.Range("I27").Value = 745
$ws.Range("J27").Value = 1663.1
$ws.Range("K27").Value = 745
$ws.Range("L27").Value = 1663.1
$ws.Range("M27").Value = -638
$ws.Range("N27").Value = -1877.1
$ws.Range("H61").Value = 2079.7778
$ws.Range("I61").Value = 1601.6
$ws.Range("J61").Value = 2677.5
$ws.Range("K61").Value = 1601.6
$ws.Range("L61").Value = 2677.5
$ws.Range("M61").Value = -1399.6
$ws.Range("N61").Value = -3081.5
$ws.Range("H113").Value = 2079.7778
$ws.Range("I113").Value = 1601.6
$ws.Range("J113").Value = 2677.5
$ws.Range("K113").Value = 1601.6
$ws.Range("L113").Value = 2677.5
$ws.Range("M113").Value = 568.4000000000001
$ws.Range("N113").Value = -7017.5
$ws.Range("H134").Value = 23479.666
$ws.Range("J134").Value = 23479.666
$ws.Range("L134").Value = 23479.666
$ws.Range("N134").Value = -33619.666

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166668640
$ws.Range("I62").Value = 250001460
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 250001460
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -250000836
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 166668640
$ws.Range("I65").Value = 250001460
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 1250007300
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -1250004180
$ws.Range("N65").Value = -21240
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 69445680
$ws.Range("I126").Value = 85471250
$ws.Range("K126").Value = 256413750
$ws.Range("M126").Value = -256411280
$ws.Range("H132").Value = 1656
$ws.Range("I132").Value = 1021.3
$ws.Range("K132").Value = 3063.9
$ws.Range("M132").Value = -533.8999999999996
